$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# 1) Convert the "git clone <url>" w:hyperlink into a field-code hyperlink
#    (fldChar begin / instrText HYPERLINK "..." / fldChar separate / display
#    text / fldChar end), keeping the same external relationship/target.
# ---------------------------------------------------------------------------
$hl = $d.Hyperlinks.Item(1)
$url = $hl.Address
$rng = $hl.Range

$frag = '<w:r xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:fldChar w:fldCharType="begin"/></w:r>' +
        '<w:r xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:rPr><w:lang w:val="en-US"/></w:rPr><w:instrText xml:space="preserve"> HYPERLINK "' + $url + '" </w:instrText></w:r>' +
        '<w:r xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:fldChar w:fldCharType="separate"/></w:r>' +
        '<w:r xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:rPr><w:rStyle w:val="Hyperlink"/><w:lang w:val="en-US"/></w:rPr><w:t>' + $url + '</w:t></w:r>' +
        '<w:r xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:rPr><w:rStyle w:val="Hyperlink"/><w:lang w:val="en-US"/></w:rPr><w:fldChar w:fldCharType="end"/></w:r>'

$xml = '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">' +
       '<pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">' +
       '<pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p>' + $frag + '</w:p></w:body></w:document></pkg:xmlData>' +
       '</pkg:part></pkg:package>'

$rng.InsertXML($xml)

# ---------------------------------------------------------------------------
# 2) Move the "_GoBack" bookmark from the end of the document to right after
#    the second "DTure" (immediately before the following proofErr/spellEnd),
#    and drop the single-space run that used to follow it there.
# ---------------------------------------------------------------------------
$goBack = $d.Bookmarks.Item("_GoBack")
$goBack.Delete()

$full = $d.Range(0, $d.Content.End)
$full.Find.Execute("DTure", $false, $false, $false, $false, $false, $true, 1, $false, "", 0)
$second = $d.Range($full.End, $d.Content.End)
$second.Find.Execute("DTure", $false, $false, $false, $false, $false, $true, 1, $false, "", 0)

$spaceRng = $d.Range($second.End, $second.End + 1)
$spaceRng.Delete()

$newBmRange = $d.Range($second.End, $second.End)
$d.Bookmarks.Add("_GoBack", $newBmRange)
